$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 5

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 2

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 4
